$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 214, shifting existing rows 214-329 down to 215-330.
$ws.Rows.Item(214).Insert()

# Populate the newly inserted row 214 with the new record's data.
$ws.Range("A214").Value = 10
$ws.Range("B214").Value = "Vega Modelo de Temuco"
$ws.Range("C214").Value = "La Araucanía"
$ws.Range("D214").Value = 45029
$ws.Range("E214").Value = 9
$ws.Range("F214").Value = 100112043
$ws.Range("G214").Value = "Pepino dulce"
$ws.Range("H214").Value = "Cultivar IV Región"
$ws.Range("I214").Value = "Primera"
$ws.Range("J214").Value = 330
$ws.Range("K214").Value = 16000
$ws.Range("L214").Value = 17000
$ws.Range("M214").Value = 16545
$ws.Range("N214").Value = "$/bandeja 18 kilos"
$ws.Range("O214").Value = "Provincia de Limarí"
$ws.Range("P214").Value = 919
$ws.Range("Q214").Value = 18
$ws.Range("R214").Value = "Hortaliza"
